# September_Schedule.xlsx edit
# - Trim the "Benji's defense..." note (drop "if time allows")
# - Fill in week-of-9/9 (row 4) Morning/Afternoon cells with the qPCR /
#   extractions plan, matching row 3's wording for the same columns
# - Resize rows 3 & 4 to their new wrapped-text heights
# - Leave the cursor parked on L7, like the author's last save

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rightQuote = [char]0x2019

# Shorten the shared string used by G3 ("Benji's defense, ...").
$ws.Range("G3").Value2 = "Benji" + $rightQuote + "s defense, derivatisation with Zsolt"

# Populate row 4 (week beginning 09/09) the same way row 3 is laid out:
# Morning/Afternoon for Mon-Thu = "qPCR", Friday Morning = "Extractions and
# concentrations", Friday Afternoon = "Extractions and concentrations, data
# analysis".
$ws.Range("B4").Value2 = "qPCR"
$ws.Range("C4").Value2 = "qPCR"
$ws.Range("D4").Value2 = "qPCR"
$ws.Range("E4").Value2 = "qPCR"
$ws.Range("F4").Value2 = "qPCR"
$ws.Range("G4").Value2 = "qPCR"
$ws.Range("H4").Value2 = "qPCR"
$ws.Range("I4").Value2 = "qPCR"
$ws.Range("J4").Value2 = "Extractions and concentrations"
$ws.Range("K4").Value2 = "Extractions and concentrations, data analysis"

# Rows re-wrap to their new content heights.
$ws.Rows.Item(3).RowHeight = 102.2
$ws.Rows.Item(4).RowHeight = 57.45

# Restore the author's final selection.
$ws.Range("L7").Select() | Out-Null
